$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells: "_old" -> "_FV2310" and "_new" -> "_FV2404"
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Freeze the header row
$ws.Activate()
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true

# Create a table (ListObject) over the used range
$rng = $ws.Range("A1:U60")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
